$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.047805713851277
$ws.Range("D2").Value = 1.055831166408274
$ws.Range("E2").Value = 1.055693980354425
$ws.Range("F2").Value = 1.067563987234846
$ws.Range("I2").Value = 1.046581727576755
$ws.Range("J2").Value = 1.052852826138828
$ws.Range("K2").Value = 1.058570543986037
$ws.Range("L2").Value = 1.058433734901368
$ws.Range("M2").Value = 1.070271512438324
$ws.Range("N2").Value = 1.054347997915852
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.048660813739666
$ws.Range("D3").Value = 1.056510869454215
$ws.Range("E3").Value = 1.05644292037201
$ws.Range("F3").Value = 1.068365478828986
$ws.Range("I3").Value = 1.046793066235077
$ws.Range("J3").Value = 1.053357011502135
$ws.Range("K3").Value = 1.05906410726893
$ws.Range("L3").Value = 1.058996331465311
$ws.Range("M3").Value = 1.070888849259337
$ws.Range("N3").Value = 1.054852899280207
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.049214869018332
$ws.Range("D4").Value = 1.056951346823777
$ws.Range("E4").Value = 1.05692855212181
$ws.Range("F4").Value = 1.06888515607793
$ws.Range("I4").Value = 1.04692904354447
$ws.Range("J4").Value = 1.053683321878806
$ws.Range("K4").Value = 1.059383457344274
$ws.Range("L4").Value = 1.05936071785764
$ws.Range("M4").Value = 1.071288706662163
$ws.Range("N4").Value = 1.05517967305504
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.049447971343948
$ws.Range("D5").Value = 1.057136680645537
$ws.Range("E5").Value = 1.057132953398287
$ws.Range("F5").Value = 1.069103879819791
$ws.Range("I5").Value = 1.046986022605227
$ws.Range("J5").Value = 1.053820517884157
$ws.Range("K5").Value = 1.059517706332912
$ws.Range("L5").Value = 1.059513987913921
$ws.Range("M5").Value = 1.071456900645536
$ws.Range("N5").Value = 1.055317063894456
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.049487120608688
$ws.Range("D6").Value = 1.057167808170738
$ws.Range("E6").Value = 1.057167287389581
$ws.Range("F6").Value = 1.069140619184002
$ws.Range("I6").Value = 1.046995578722844
$ws.Range("J6").Value = 1.053843554545265
$ws.Range("K6").Value = 1.059540246949462
$ws.Range("L6").Value = 1.059539727397117
$ws.Range("M6").Value = 1.071485146615987
$ws.Range("N6").Value = 1.055340133270264
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.049217983047738
$ws.Range("D7").Value = 1.056953822648251
$ws.Range("E7").Value = 1.056931282394383
$ws.Range("F7").Value = 1.068888077689856
$ws.Range("I7").Value = 1.04692980563248
$ws.Range("J7").Value = 1.05368515504079
$ws.Range("K7").Value = 1.05938525121029
$ws.Range("L7").Value = 1.059362765539228
$ws.Range("M7").Value = 1.071290953711372
$ws.Range("N7").Value = 1.055181508820325
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.048094542863258
$ws.Range("D8").Value = 1.05606073663402
$ws.Range("E8").Value = 1.055946876323669
$ws.Range("F8").Value = 1.067834634678712
$ws.Range("I8").Value = 1.0466533100506
$ws.Range("J8").Value = 1.053023202910428
$ws.Range("K8").Value = 1.058737348979744
$ws.Range("L8").Value = 1.058623793930598
$ws.Range("M8").Value = 1.07048006066557
$ws.Range("N8").Value = 1.05451861664201
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.046120705406469
$ws.Range("D9").Value = 1.054492172116043
$ws.Range("E9").Value = 1.054220103648379
$ws.Range("F9").Value = 1.065986529672524
$ws.Range("I9").Value = 1.046160203362361
$ws.Range("J9").Value = 1.051857344030366
$ws.Range("K9").Value = 1.057595583036072
$ws.Range("L9").Value = 1.057324368619307
$ws.Range("M9").Value = 1.069054293379266
$ws.Range("N9").Value = 1.053351102108613
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.044808819968652
$ws.Range("D10").Value = 1.053450046420203
$ws.Range("E10").Value = 1.053074331041592
$ws.Range("F10").Value = 1.064760090044467
$ws.Range("I10").Value = 1.045827553446746
$ws.Range("J10").Value = 1.051080577883082
$ws.Range("K10").Value = 1.056834435079203
$ws.Range("L10").Value = 1.056460014086817
$ws.Range("M10").Value = 1.068105984092427
$ws.Range("N10").Value = 1.052573232864303
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.044241730678205
$ws.Range("D11").Value = 1.05299966782085
$ws.Range("E11").Value = 1.052579505844349
$ws.Range("F11").Value = 1.064230388941324
$ws.Range("I11").Value = 1.045682594039424
$ws.Range("J11").Value = 1.050744358582581
$ws.Range("K11").Value = 1.05650487405885
$ws.Range("L11").Value = 1.056086215332801
$ws.Range("M11").Value = 1.067695898748992
$ws.Range("N11").Value = 1.05223653609383
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.044031235279619
$ws.Range("D12").Value = 1.052832509557397
$ws.Range("E12").Value = 1.052395903283041
$ws.Range("F12").Value = 1.064033840053496
$ws.Range("I12").Value = 1.045628612232929
$ws.Range("J12").Value = 1.050619491983034
$ws.Range("K12").Value = 1.056382464997327
$ws.Range("L12").Value = 1.055947442332586
$ws.Range("M12").Value = 1.067543657423469
$ws.Range("N12").Value = 1.052111492169389
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.044076380585898
$ws.Range("D13").Value = 1.052868359550588
$ws.Range("E13").Value = 1.052435277684843
$ws.Range("F13").Value = 1.064075991116293
$ws.Range("I13").Value = 1.045640197720346
$ws.Range("J13").Value = 1.050646275355818
$ws.Range("K13").Value = 1.056408721926873
$ws.Range("L13").Value = 1.055977206307036
$ws.Range("M13").Value = 1.067576309938493
$ws.Range("N13").Value = 1.052138313577635
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.044224328057985
$ws.Range("D14").Value = 1.052985847754781
$ws.Range("E14").Value = 1.052564325158126
$ws.Range("F14").Value = 1.064214137941077
$ws.Range("I14").Value = 1.045678134685545
$ws.Range("J14").Value = 1.050734036650219
$ws.Range("K14").Value = 1.05649475558615
$ws.Range("L14").Value = 1.056074742825546
$ws.Range("M14").Value = 1.06768331273821
$ws.Range("N14").Value = 1.052226199503141
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.044315502884127
$ws.Range("D15").Value = 1.053058253715927
$ws.Range("E15").Value = 1.052643861801077
$ws.Range("F15").Value = 1.064299282076381
$ws.Range("I15").Value = 1.045701490714103
$ws.Range("J15").Value = 1.050788111999888
$ws.Range("K15").Value = 1.056547764411535
$ws.Range("L15").Value = 1.056134847946134
$ws.Range("M15").Value = 1.067749251694077
$ws.Range("N15").Value = 1.052280351646009
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.044846476311298
$ws.Range("D16").Value = 1.053479955028383
$ws.Range("E16").Value = 1.053107198556583
$ws.Range("F16").Value = 1.064795273293741
$ws.Range("I16").Value = 1.045837154611008
$ws.Range("J16").Value = 1.051102894411839
$ws.Range("K16").Value = 1.056856307522955
$ws.Range("L16").Value = 1.056484831962885
$ws.Range("M16").Value = 1.068133211601807
$ws.Range("N16").Value = 1.052595581085091
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.045179801687818
$ws.Range("D17").Value = 1.053744711030778
$ws.Range("E17").Value = 1.053398187397389
$ws.Range("F17").Value = 1.065106759801028
$ws.Range("I17").Value = 1.045922007269658
$ws.Range("J17").Value = 1.05130038367572
$ws.Range("K17").Value = 1.057049855100739
$ws.Range("L17").Value = 1.056704495253553
$ws.Range("M17").Value = 1.068374204952655
$ws.Range("N17").Value = 1.05279335080638
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.045374317806987
$ws.Range("D18").Value = 1.053899222418414
$ws.Range("E18").Value = 1.053568041748896
$ws.Range("F18").Value = 1.065288575244153
$ws.Range("I18").Value = 1.045971411591456
$ws.Range("J18").Value = 1.051415587846663
$ws.Range("K18").Value = 1.057162750066442
$ws.Range("L18").Value = 1.056832666618545
$ws.Range("M18").Value = 1.068514824174337
$ws.Range("N18").Value = 1.052908718580458
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.045440658510331
$ws.Range("D19").Value = 1.053951920946596
$ws.Range("E19").Value = 1.053625978918126
$ws.Range("F19").Value = 1.065350591706698
$ws.Range("I19").Value = 1.045988242103596
$ws.Range("J19").Value = 1.051454871490213
$ws.Range("K19").Value = 1.057201244618135
$ws.Range("L19").Value = 1.056876377410437
$ws.Range("M19").Value = 1.068562780443322
$ws.Range("N19").Value = 1.052948058011288
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.045144029390778
$ws.Range("D20").Value = 1.053716296547199
$ws.Range("E20").Value = 1.053366954055197
$ws.Range("F20").Value = 1.065073326719791
$ws.Range("I20").Value = 1.045912912561773
$ws.Range("J20").Value = 1.051279193693203
$ws.Range("K20").Value = 1.057029089066718
$ws.Range("L20").Value = 1.056680922745768
$ws.Range("M20").Value = 1.068348343268834
$ws.Range("N20").Value = 1.052772130731656
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.044180757148353
$ws.Range("D21").Value = 1.052951246736876
$ws.Range("E21").Value = 1.05252631842467
$ws.Range("F21").Value = 1.064173451449888
$ws.Range("I21").Value = 1.045666966982419
$ws.Range("J21").Value = 1.050708192567396
$ws.Range("K21").Value = 1.056469420680204
$ws.Range("L21").Value = 1.056046018738319
$ws.Range("M21").Value = 1.067651800786213
$ws.Range("N21").Value = 1.052200318718756
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.043575959929088
$ws.Range("D22").Value = 1.052470996524158
$ws.Range("E22").Value = 1.051998921362597
$ws.Range("F22").Value = 1.063608855211856
$ws.Range("I22").Value = 1.045511536495727
$ws.Range("J22").Value = 1.050349299147233
$ws.Range("K22").Value = 1.056117562121106
$ws.Range("L22").Value = 1.05564724960613
$ws.Range("M22").Value = 1.067214335783532
$ws.Range("N22").Value = 1.05184091562877
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.04389649320719
$ws.Range("D23").Value = 1.052725512918637
$ws.Range("E23").Value = 1.05227839544138
$ws.Range("F23").Value = 1.063908044773441
$ws.Range("I23").Value = 1.045594008214959
$ws.Range("J23").Value = 1.050539543735568
$ws.Range("K23").Value = 1.056304085954447
$ws.Range("L23").Value = 1.055858604329428
$ws.Range("M23").Value = 1.067446198304523
$ws.Range("N23").Value = 1.052031430386243
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.045160193064661
$ws.Range("D24").Value = 1.053729135570899
$ws.Range("E24").Value = 1.053381066670644
$ws.Range("F24").Value = 1.065088433286728
$ws.Range("I24").Value = 1.045917022343257
$ws.Range("J24").Value = 1.051288768497209
$ws.Range("K24").Value = 1.05703847233782
$ws.Range("L24").Value = 1.056691574007028
$ws.Range("M24").Value = 1.068360028890474
$ws.Range("N24").Value = 1.052781719132982
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.046630291023701
$ws.Range("D25").Value = 1.054897059982211
$ws.Range("E25").Value = 1.05466557067015
$ws.Range("F25").Value = 1.066463325487791
$ws.Range("I25").Value = 1.046288376131126
$ws.Range("J25").Value = 1.052158668935996
$ws.Range("K25").Value = 1.057890757790195
$ws.Range("L25").Value = 1.057659967706949
$ws.Range("M25").Value = 1.069422507407023
$ws.Range("N25").Value = 1.053652854930171
